# Sprint 4 Backlog - Burndown: mark the "Complete functionality to remove
# ingredients used to cook recipe from pantry (desktop)" task (row 14) as
# having 1.5 hrs of completed time logged against Janera, and reassign the
# "Add new recipe" task in row 16 from Matthew to Janera.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: log 1.5 hrs completed, attributed to Janera.
$ws.Range("E14").Value = 1.5
$ws.Range("F14").Value = "Janera"

# Row 16: switch the assigned team member from Matthew to Janera.
$ws.Range("D16").Value = "Janera"

# Reflect the edit in the sheet's active selection, as the author left it.
$ws.Activate()
$ws.Range("F16").Select()
